$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: "Er verstarb am STERBEDATUM in STERBEORT. " ->
#         "PERSONALPRONOMEN" + "N" + " verstarb am STERBEDATUM in STERBEORT. "
# The run containing "Er" is split in three: the two-letter "Er" is
# replaced by two new runs ("PERSONALPRONOMEN" and "N"); formatting is
# round-tripped (Bold on/off) so the engine keeps the runs distinct
# without leaving any residual property behind.
# -----------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute("Er verstarb am STERBEDATUM in STERBEORT. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1 = $find1.Parent
$erRange = $d.Range($rng1.Start, $rng1.Start + 2)
$erRange.Text = "PERSONALPRONOMEN"
$erRange.Font.Bold = 1
$erRange.Font.Bold = 0

$nRange = $d.Range($erRange.End, $erRange.End)
$nRange.InsertAfter("N")
$nRange2 = $d.Range($erRange.End, $erRange.End + 1)
$nRange2.Font.Bold = 1
$nRange2.Font.Bold = 0

# -----------------------------------------------------------------
# Edit 2: "hat uns auch" -> "hat, uns auch"
# A comma is inserted as its own run between "hat" and " uns auch...".
# -----------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("zugesagt hat uns auch", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2 = $find2.Parent
$hatEnd = $rng2.Start + ("zugesagt hat").Length
$commaRange = $d.Range($hatEnd, $hatEnd)
$commaRange.InsertAfter(",")
$commaRange2 = $d.Range($hatEnd, $hatEnd + 1)
$commaRange2.Font.Bold = 1
$commaRange2.Font.Bold = 0

Write-Host "Final text:" $d.Content.Text
